# Background Research on Quantum Algorithms
# Adds two new reference rows (14 & 15) to the References sheet:
#   Row 14: Grover's Algorithm YouTube explanation (hyperlinked) + description
#   Row 15: Qiskit Textbook citation (plain text) + long BibTeX-style citation
#
# Also updates the active selection to reflect where the author was last
# working (cell A13).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 14 - Grover's Algorithm video
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("A14"), "https://www.youtube.com/watch?v=9Iorl3MUF-Y")
$ws.Range("A14").Value = "https://www.youtube.com/watch?v=9Iorl3MUF-Y"
# Hyperlinks.Add mints its own (wrapped) style; restore the plain Hyperlink
# style used by every other linked cell in column A.
$ws.Range("A14").Style = "Hyperlink"

$ws.Range("B14").Value = "Grovers Algorithm, Best explanation by Umesh Vazirani"

$ws.Rows.Item(14).RowHeight = 43.2

# ---------------------------------------------------------------------------
# Row 15 - Qiskit Textbook citation
# ---------------------------------------------------------------------------
$ws.Range("A15").Value = "Qiskit Textbook"

$qiskitCitation = @'
Qiskit-Textbook,
       author = {Abraham Asfaw, Luciano Bello, Yael Ben-Haim, Sergey Bravyi, Nicholas Bronn, Lauren Capelluto, Almudena Carrera Vazquez, Jack Ceroni,  Richard Chen, Albert Frisch, Jay Gambetta, Shelly Garion, Leron Gil, Salvador De La Puente Gonzalez, Francis Harkins, Takashi Imamichi, David McKay, Antonio Mezzacapo, Zlatko Minev, Ramis Movassagh, Giacomo Nannicni, Paul Nation,  Anna Phan, Marco Pistoia, Arthur Rattew, Joachim Schaefer, Javad Shabani, John Smolin, Kristan Temme, Madeleine Tod, Stephen Wood, James Wootton.},
       title = {Learn Quantum Computation Using Qiskit},
       year = {2020},
       url = {http://community.qiskit.org/textbook},
@misc{ Qiskit-Textbook,
       author = {Abraham Asfaw, Luciano Bello, Yael Ben-Haim, Sergey Bravyi, Nicholas Bronn, Lauren Capelluto, Almudena Carrera Vazquez, Jack Ceroni,  Richard Chen, Albert Frisch, Jay Gambetta, Shelly Garion, Leron Gil, Salvador De La Puente Gonzalez, Francis Harkins, Takashi Imamichi, David McKay, Antonio Mezzacapo, Zlatko Minev, Ramis Movassagh, Giacomo Nannicni, Paul Nation,  Anna Phan, Marco Pistoia, Arthur Rattew, Joachim Schaefer, Javad Shabani, John Smolin, Kristan Temme, Madeleine Tod, Stephen Wood, James Wootton.},
       title = {Learn Quantum Computation Using Qiskit},
       year = {2020},
       url = {http://community.qiskit.org/textbook},
@misc{ Qiskit-Textbook,
       author = {Abraham Asfaw, Luciano Bello, Yael Ben-Haim, Sergey Bravyi, Nicholas Bronn, Lauren Capelluto, Almudena Carrera Vazquez, Jack Ceroni,  Richard Chen, Albert Frisch, Jay Gambetta, Shelly Garion, Leron Gil, Salvador De La Puente Gonzalez, Francis Harkins, Takashi Imamichi, David McKay, Antonio Mezzacapo, Zlatko Minev, Ramis Movassagh, Giacomo Nannicni, Paul Nation,  Anna Phan, Marco Pistoia, Arthur Rattew, Joachim Schaefer, Javad Shabani, John Smolin, Kristan Temme, Madeleine Tod, Stephen Wood, James Wootton.},
       title = {Learn Quantum Computation Using Qiskit},
       year = {2020},
       url = {http://community.qiskit.org/textbook},
}
'@

$ws.Range("B15").Value = $qiskitCitation

$ws.Rows.Item(15).RowHeight = 48.6

# ---------------------------------------------------------------------------
# Selection - leave the active cell where the author last left it
# ---------------------------------------------------------------------------
$ws.Range("A13").Select()
